$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 3.1
$ws.Range("J3").Value = 3
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 6.5
$ws.Range("AC3").Value = 8
$ws.Range("AE3").Value = 17
$ws.Range("AM3").Value = 351
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 26
$ws.Range("AT3").Value = 2.5
$ws.Range("AV3").Value = 67
$ws.Range("BA3").Value = 101
